# Add 4 new data rows to Table1 on the "patchtables" sheet:
#   - subh    (inserted right after "subv",   before "subv25")
#   - subh30  (inserted right after "subv30", before "subv40")
#   - comh36  (appended after "comv36")
#   - subv27  (appended after "comh36")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

function Fill-Row {
    param($row, $spec, $extra, $orient, $width, $color)

    $ws.Cells.Item($row, 1).Formula = "=Table1[[#This Row],[spec]]&E$row&C$row"
    $ws.Cells.Item($row, 2).Value = $spec
    if ($null -eq $extra) {
        $ws.Cells.Item($row, 3).ClearContents()
    } else {
        $ws.Cells.Item($row, 3).Value = $extra
    }
    $ws.Cells.Item($row, 4).Formula = '=IF(Table1[[#This Row],[spec]]="com","36/72 (6x12)",IF(Table1[[#This Row],[spec]]="sub","27/60 (6x10)","27/60 (5x12)"))&IF(Table1[[#This Row],[hooktype]]="",""," +"&Table1[[#This Row],[hooktype]])'
    $ws.Cells.Item($row, 5).Value = $orient
    if ($null -eq $width) {
        $ws.Cells.Item($row, 6).ClearContents()
    } else {
        $ws.Cells.Item($row, 6).Value = $width
    }
    $ws.Cells.Item($row, 7).Formula = "=IF(E$row=""v"", VLOOKUP(B$row, `$K`$2:`$M`$4,2)+F$row, VLOOKUP(B$row,`$K`$2:`$M`$4,3))"
    $ws.Cells.Item($row, 8).Formula = "=IF(E$row=""v"", VLOOKUP(B$row,`$K`$2:`$M`$4, 3), VLOOKUP(B$row,`$K`$2:`$M`$4,2)+F$row)"
    $ws.Cells.Item($row, 9).Value = $color
}

# --- Step 1: insert "subh" as the new row 3 (between "subv" and "subv25") ---
$lo.Resize($ws.Range("A1:I11"))
$ws.Range("A3:I3").Insert()
# The row-insert shifted the K:M helper lookup table (rows 2-4) down by one row too;
# that table must stay put at K2:M4, so shift it back up and clear the spilled-over row.
for ($c = 11; $c -le 13; $c++) {
    $ws.Cells.Item(3, $c).Value = $ws.Cells.Item(4, $c).Value()
    $ws.Cells.Item(4, $c).Value = $ws.Cells.Item(5, $c).Value()
}
$ws.Range("K5:M5").ClearContents()
Fill-Row 3 "sub" $null "h" $null "lightgreen"

# --- Step 2: insert "subh30" as the new row 6 (between "subv30" and "subv40") ---
$lo.Resize($ws.Range("A1:I12"))
$ws.Range("A6:I6").Insert()
Fill-Row 6 "sub" 30 "h" 36 "lime"

# --- Step 3: append "comh36" as the new last row (13) ---
$lo.Resize($ws.Range("A1:I13"))
Fill-Row 13 "com" 36 "h" 36 "purple"

# --- Step 4: append "subv27" as the new last row (14) ---
$lo.Resize($ws.Range("A1:I14"))
Fill-Row 14 "sub" 27 "v" 36 "springgreen"

$ws.Range("E14").Select()
